# Applies updated crypto price/volume figures, and restores the two
# rows whose coin order was swapped (Algorand/MultiversX and
# TrustWalletToken/NEARProtocol), per the Dec 16 2023 data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.267.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.22%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.248.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.35%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.56%  "

$ws.Range("E6").Value = "  -2.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.60"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.37%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  -4.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.16"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.16%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0946"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.07%  "

$ws.Range("E12").Value = "  -2.85%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.102"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.73%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.854"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.78%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.250.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.98%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "42.127.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.36%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0979"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.52%  "

$ws.Range("E21").Value = "  +4.60%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "232.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.96%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +38.26%  "

$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.52%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.05%  "

$ws.Range("E27").Value = "  -3.73%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "169.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.22%  "

$ws.Range("E29").Value = "  -4.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0826"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.72%  "

$ws.Range("E32").Value = "  -5.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.33"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.13%  "

$ws.Range("E34").Value = "  -1.88%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.17"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.32%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0306"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "13.47"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.45%  "

$ws.Range("E39").Value = "  -4.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.66%  "

$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "61.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.11%  "

$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.203"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "108.98"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.05%  "

$ws.Range("E44").Value = "  -3.33%  "

$ws.Range("E45").Value = "  +0.46%  "

$ws.Range("E46").Value = "  -0.45%  "

$ws.Range("E47").Value = "  -3.40%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.97%  "

$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.89%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.13"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.96%  "

$ws.Range("E51").Value = "  +0.46%  "
